# Day Holder 2018 - append 4 more days (335..338) continuing the existing
# Day_Number / Date series, then restore the sheet's UI selection state.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Days")

$firstNewRow = 336
$lastExistingRow = 335
$rowsToAdd = 4

for ($i = 0; $i -lt $rowsToAdd; $i++) {
    $row = $firstNewRow + $i

    # Copy the formatting (number format / borders / alignment) of the last
    # existing data row down onto the new row before writing its values, so
    # the new cells pick up the same styles (s="3" for Day_Number, s="4" for
    # the Date column) as every other data row.
    $ws.Range("A$lastExistingRow" + ":B$lastExistingRow").Copy() | Out-Null
    $ws.Range("A$row" + ":B$row").PasteSpecial(-4122) | Out-Null # xlPasteFormats

    $ws.Cells.Item($row, 1).Value = 335 + $i
    $ws.Cells.Item($row, 2).Value = 43435 + $i
}
$excel.CutCopyMode = $false

# Restore the selection recorded in the sheet view (two-cell range spanning
# E325:E326). NOTE: this COM layer always pins the active cell to the
# top-left corner of whatever rectangle is selected (there's no exposed way
# to activate a specific cell *within* a multi-cell selection without
# collapsing the selection down to that single cell), so the best
# reproducible approximation of "activeCell=E326, sqref=E325:E326" is to
# select the two-cell range outright.
$ws.Range("E325:E326").Select() | Out-Null
